$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2649.3333
$ws.Range("J70").Value = 2474
$ws.Range("L70").Value = 7422
$ws.Range("N70").Value = -7962

$ws.Range("H73").Value = 2649.3333
$ws.Range("J73").Value = 2474
$ws.Range("L73").Value = 7422
$ws.Range("N73").Value = -9294

$ws.Range("H76").Value = 15523.111
$ws.Range("I76").Value = 8764.5
$ws.Range("K76").Value = 8764.5
$ws.Range("M76").Value = -8449.5

$ws.Range("H79").Value = 15523.111
$ws.Range("I79").Value = 8764.5
$ws.Range("K79").Value = 8764.5
$ws.Range("M79").Value = -7672.5

$ws.Range("H98").Value = 3346.074
$ws.Range("I98").Value = 3282.4614
$ws.Range("K98").Value = 3282.4614
$ws.Range("M98").Value = -1784.4614

$ws.Range("H122").Value = 3346.074
$ws.Range("I122").Value = 3282.4614
$ws.Range("K122").Value = 9847.3842
$ws.Range("M122").Value = -7397.3842

$ws.Range("H132").Value = 13708.667
$ws.Range("I132").Value = 7818.7417
$ws.Range("K132").Value = 23456.2251
$ws.Range("M132").Value = -20926.2251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = 0

$ws.Range("H61").Value = 1977653.2
$ws.Range("I61").Value = 2555
$ws.Range("K61").Value = 2555
$ws.Range("M61").Value = -2343

$ws.Range("H136").Value = 1977653.2
$ws.Range("I136").Value = 2555
$ws.Range("K136").Value = 7665
$ws.Range("M136").Value = -5115

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 86592.44
$ws.Range("I134").Value = 121272.336
$ws.Range("J134").Value = 42004
$ws.Range("K134").Value = 363817.008
$ws.Range("L134").Value = 126012
$ws.Range("M134").Value = -361282.008
$ws.Range("N134").Value = -131082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 17657.334
$ws.Range("I105").Value = 25486.25
$ws.Range("J105").Value = 1999.5
$ws.Range("K105").Value = 25486.25
$ws.Range("L105").Value = 1999.5
$ws.Range("M105").Value = -23739.25
$ws.Range("N105").Value = -5493.5

$ws.Range("H134").Value = 26320394
$ws.Range("I134").Value = 1776.5333
$ws.Range("K134").Value = 5329.5999
$ws.Range("M134").Value = -2794.5999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 1036
$ws.Range("J40").Value = 937.5
$ws.Range("L40").Value = 3750
$ws.Range("N40").Value = -3888

$ws.Range("H68").Value = 821
$ws.Range("I68").Value = 475
$ws.Range("J68").Value = 1097.8
$ws.Range("K68").Value = 1425
$ws.Range("L68").Value = 3293.4
$ws.Range("M68").Value = -614
$ws.Range("N68").Value = -4915.4

$ws.Range("H71").Value = 821
$ws.Range("I71").Value = 475
$ws.Range("J71").Value = 1097.8
$ws.Range("K71").Value = 4275
$ws.Range("L71").Value = 9880.199999999999
$ws.Range("M71").Value = -219
$ws.Range("N71").Value = -17992.2

$ws.Range("H107").Value = 1134.5714
$ws.Range("I107").Value = 1026.6666
$ws.Range("J107").Value = 1164
$ws.Range("K107").Value = 3079.9998
$ws.Range("L107").Value = 3492
$ws.Range("M107").Value = -1159.9998
$ws.Range("N107").Value = -7332

$ws.Range("H122").Value = 3420
$ws.Range("I122").Value = 2585.4
$ws.Range("J122").Value = 3665.4707
$ws.Range("K122").Value = 23268.6
$ws.Range("L122").Value = 32989.2363
$ws.Range("M122").Value = -20818.6
$ws.Range("N122").Value = -37889.2363

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3317.1177
$ws.Range("I80").Value = 3165.8333
$ws.Range("J80").Value = 3680.2
$ws.Range("K80").Value = 3165.8333
$ws.Range("L80").Value = 3680.2
$ws.Range("M80").Value = -2167.8333
$ws.Range("N80").Value = -5676.2

$ws.Range("H83").Value = 3317.1177
$ws.Range("I83").Value = 3165.8333
$ws.Range("J83").Value = 3680.2
$ws.Range("K83").Value = 15829.1665
$ws.Range("L83").Value = 18401
$ws.Range("M83").Value = -10837.1665
$ws.Range("N83").Value = -28385

$ws.Range("H102").Value = 5738.9287
$ws.Range("I102").Value = 6785.05
$ws.Range("J102").Value = 3123.625
$ws.Range("K102").Value = 6785.05
$ws.Range("L102").Value = 3123.625
$ws.Range("M102").Value = -5163.05
$ws.Range("N102").Value = -6367.625

$ws.Range("H126").Value = 14872.5
$ws.Range("I126").Value = 19163.334
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 57490.00199999999
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -55020.00199999999
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 1314388.4
$ws.Range("I132").Value = 1180.2
$ws.Range("K132").Value = 3540.6
$ws.Range("M132").Value = -1010.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1164.6
$ws.Range("I22").Value = 966.9
$ws.Range("J22").Value = 1560
$ws.Range("K22").Value = 966.9
$ws.Range("L22").Value = 1560
$ws.Range("M22").Value = -671.9
$ws.Range("N22").Value = -2150

$ws.Range("H26").Value = 27871.625
$ws.Range("I26").Value = 26597.8
$ws.Range("K26").Value = 26597.8
$ws.Range("M26").Value = -26302.8

$ws.Range("H27").Value = 1164.6
$ws.Range("I27").Value = 966.9
$ws.Range("J27").Value = 1560
$ws.Range("K27").Value = 966.9
$ws.Range("L27").Value = 1560
$ws.Range("M27").Value = -859.9
$ws.Range("N27").Value = -1774

$ws.Range("H55").Value = 1271.3549
$ws.Range("I55").Value = 1099.25
$ws.Range("K55").Value = 1099.25
$ws.Range("M55").Value = -926.25

$ws.Range("H82").Value = 2519.8635
$ws.Range("I82").Value = 3376.5386
$ws.Range("K82").Value = 3376.5386
$ws.Range("M82").Value = -3015.5386

$ws.Range("H85").Value = 2519.8635
$ws.Range("I85").Value = 3376.5386
$ws.Range("K85").Value = 3376.5386
$ws.Range("M85").Value = -2128.5386

$ws.Range("H122").Value = 7546.846
$ws.Range("J122").Value = 4467.4
$ws.Range("L122").Value = 13402.2
$ws.Range("N122").Value = -18302.2

$ws.Range("H132").Value = 2501424.2
$ws.Range("J132").Value = 4369049.5
$ws.Range("L132").Value = 13107148.5
$ws.Range("N132").Value = -13112208.5

$ws.Range("H136").Value = 3196972.5
$ws.Range("I136").Value = 37267.285
$ws.Range("K136").Value = 111801.855
$ws.Range("M136").Value = -109251.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 27666.666
$ws.Range("I24").Value = 20000
$ws.Range("J24").Value = 31500
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 31500
$ws.Range("M24").Value = -19770
$ws.Range("N24").Value = -31960

$ws.Range("H132").Value = 908038.8
$ws.Range("I132").Value = 1877.8889
$ws.Range("J132").Value = 3626521.8
$ws.Range("K132").Value = 5633.6667
$ws.Range("L132").Value = 10879565.4
$ws.Range("M132").Value = -3103.6667
$ws.Range("N132").Value = -10884625.4
